# Potato3d FPS log update: add "Front to back" removal benchmark data block
# (cache-hit stats, wide-FOV / cache comparison table, and misc summary rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Main comparison block, rows 164-177: columns A (label), B/C/D (raw frame
# counts for three variants), E/F/G (= B/C/D divided by 30 -> seconds),
# I/J/K (scaled projections of E/F/G, only present from row 169 onward).
# ---------------------------------------------------------------------------
$mainRows = @(
    @{r=164; a="16x";             b=366; c=263; d=273}
    @{r=165; a="8x";              b=355; c=255; d=263}
    @{r=166; a="4x";              b=359; c=259; d=268}
    @{r=167; a="FOV";             b=361; c=261            }
    @{r=168                                               }
    @{r=169; a="Wide FOV";        b=295; c=238; d=247; ik=254; jk=352; kk=407}
    @{r=170; a="No tex cache";    b=312; c=260; d=268; ik=254; jk=352; kk=407}
    @{r=171; a="No Vx Cache";     b=279; c=225; d=235; ik=254; jk=352; kk=407}
    @{r=172; a="Obj3d IWRAM";     b=322; c=265; d=274; ik=254; jk=352; kk=407}
    @{r=173; a="Front to back";   b=284; c=225; d=233; ik=254; jk=352; kk=407}
    @{r=174                                               }
    @{r=175; a="No early reject"; b=310; c=238; d=229; ik=254; jk=493; kk=499}
    @{r=176;                      b=320; c=271; d=267; ik=254                }
    @{r=177;                      b=329; c=285; d=278; ik=254                }
)

foreach ($row in $mainRows) {
    $r = $row.r

    if ($row.a -ne $null) { $ws.Cells.Item($r, 1).Value = $row.a }          # A
    if ($row.b -ne $null) { $ws.Cells.Item($r, 2).Value = $row.b }          # B

    if ($r -ne 168 -and $r -ne 174) {
        $cCell = $ws.Cells.Item($r, 3)
        if ($row.c -ne $null) { $cCell.Value = $row.c }
        $cCell.NumberFormat = "#,##0.0"                                    # C style
    }

    $dCell = $ws.Cells.Item($r, 4)
    if ($row.d -ne $null) { $dCell.Value = $row.d }
    $dCell.NumberFormat = "0"                                              # D style

    $eCell = $ws.Cells.Item($r, 5)
    if ($row.b -ne $null) { $eCell.Formula = "=B$r/30" }
    $eCell.NumberFormat = "0.0"                                            # E style

    $fCell = $ws.Cells.Item($r, 6)
    if ($row.c -ne $null) { $fCell.Formula = "=C$r/30" }
    $fCell.NumberFormat = "0.0"                                            # F style

    $gCell = $ws.Cells.Item($r, 7)
    if ($row.c -ne $null) { $gCell.Formula = "=D$r/30" }
    $gCell.NumberFormat = "0.0"                                            # G style

    if ($row.ik -ne $null -or $row.jk -ne $null -or $row.kk -ne $null -or $r -eq 174) {
        $iCell = $ws.Cells.Item($r, 9)
        if ($row.ik -ne $null) { $iCell.Formula = "=E$r*$($row.ik)" }
        $iCell.NumberFormat = "#,##0"                                      # I style

        $jCell = $ws.Cells.Item($r, 10)
        if ($row.jk -ne $null) { $jCell.Formula = "=F$r*$($row.jk)" }
        if ($row.jk -ne $null -or $r -eq 174 -or $r -eq 176) { $jCell.NumberFormat = "#,##0" }

        $kCell = $ws.Cells.Item($r, 11)
        if ($row.kk -ne $null) { $kCell.Formula = "=G$r*$($row.kk)" }
        if ($row.kk -ne $null -or $r -eq 174 -or $r -eq 176 -or $r -eq 177) { $kCell.NumberFormat = "#,##0" }
    }
}

# Blank styled cell left over from the old shared-formula block (no data).
$ws.Cells.Item(163, 5).NumberFormat = "0.0"

# ---------------------------------------------------------------------------
# Vertex/texture cache hit-rate blocks, rows 178-185.
# ---------------------------------------------------------------------------
$ws.Cells.Item(178, 1).Value = "Cache "

$ws.Cells.Item(179, 1).Value = "Peeks"
$ws.Cells.Item(179, 2).Value = 357680
$ws.Cells.Item(179, 3).NumberFormat = "#,##0.0"
$ws.Cells.Item(179, 3).Value = 1156197

$ws.Cells.Item(180, 1).Value = "Misses"
$ws.Cells.Item(180, 2).Value = 17681
$ws.Cells.Item(180, 3).NumberFormat = "#,##0.0"
$ws.Cells.Item(180, 3).Value = 121609

$ws.Cells.Item(181, 2).Formula = "=1-(B180/B179)"
$ws.Cells.Item(181, 2).NumberFormat = "0.0%"
$ws.Cells.Item(181, 3).Formula = "=1-(C180/C179)"
$ws.Cells.Item(181, 3).NumberFormat = "0.0%"

$ws.Cells.Item(183, 1).Value = "Sorted Nodes"
$ws.Cells.Item(183, 2).Value = 191465
$ws.Cells.Item(183, 3).NumberFormat = "#,##0.0"
$ws.Cells.Item(183, 3).Value = 1191074

$ws.Cells.Item(184, 2).Value = 9464
$ws.Cells.Item(184, 3).NumberFormat = "#,##0.0"
$ws.Cells.Item(184, 3).Value = 118945

$ws.Cells.Item(185, 2).Formula = "=1-(B184/B183)"
$ws.Cells.Item(185, 2).NumberFormat = "0.0%"
$ws.Cells.Item(185, 3).Formula = "=1-(C184/C183)"
$ws.Cells.Item(185, 3).NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Misc totals, rows 191-197.
# ---------------------------------------------------------------------------
$ws.Cells.Item(191, 2).Value = 1577090
$ws.Cells.Item(191, 3).NumberFormat = "#,##0.0"
$ws.Cells.Item(191, 3).Value = 1182119
$ws.Cells.Item(191, 4).Formula = "=C191/B191"
$ws.Cells.Item(191, 4).NumberFormat = "0.0%"

$ws.Cells.Item(192, 2).Value = 1476188
$ws.Cells.Item(192, 3).NumberFormat = "#,##0.0"
$ws.Cells.Item(192, 3).Value = 1114233
$ws.Cells.Item(192, 4).Formula = "=C192/B192"
$ws.Cells.Item(192, 4).NumberFormat = "0.0%"

$ws.Cells.Item(194, 2).Value = 4037321
$ws.Cells.Item(194, 3).NumberFormat = "#,##0.0"
$ws.Cells.Item(194, 3).Value = 1267009
$ws.Cells.Item(194, 4).Formula = "=C194/B194"
$ws.Cells.Item(194, 4).NumberFormat = "0.0%"

$ws.Cells.Item(197, 2).Value = 4355914
$ws.Cells.Item(197, 3).NumberFormat = "#,##0.0"
$ws.Cells.Item(197, 3).Value = 2371658
$ws.Cells.Item(197, 4).Formula = "=C197/B197"
$ws.Cells.Item(197, 4).NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Cosmetics: widened columns C/N/O for the new labels/table, and move the
# visible selection down to the newly-added block.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 10.71
$ws.Columns("N").ColumnWidth = 22.86
$ws.Columns("O").ColumnWidth = 27.29

$ws.Range("E177").Select()
